$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")

# Fix typo in the Crandios instructions text (row 2): "Crandio" -> "Crandios".
$ws.Range("B2").Value = "So whenever Crandios (his name or his picture) appears on the screen, hit the spacebar as quickly as possible. Hit the spacebar whenever you see an image of Crandios or the word CRANDIOS. " + [char]10 + [char]10 + "Press the spacebar to start the task"

# Replace the Palpitoad creature (row 6) with Golett.
$ws.Range("A6").Value = "For this task, the target creature is Golett:"
$ws.Range("B6").Value = "So whenever Golett (his name or his picture) appears on the screen, hit the spacebar as quickly as possible. Hit the spacebar whenever you see an image of Golett or the word GOLETT. " + [char]10 + [char]10 + "Press the spacebar to start the task"

# Update the active selection to match the saved workbook state.
$ws.Range("B2").Select()
